# Refresh the cryptos list (prices + 1h volume deltas) as produced by the
# scheduled "Updated cryptos list ... with GitHub Actions" job, plus a
# ranking swap between BitcoinCash and InternetComputer(DFINITY) at rows 22/23.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($range, $value) {
    # Excel's Value setter auto-infers numeric-looking strings (e.g. "253.39")
    # as numbers. Force the cell to Text first so the literal string is kept,
    # then restore the default "Normal" style so no stray number-format style
    # is left attached to the cell.
    $cell = $ws.Range($range)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

# Row 22 / 23 swap in ranking (BitcoinCash now above InternetComputer(DFINITY))
$ws.Range("B22").Value = "BitcoinCash"
$ws.Range("C22").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("B23").Value = "InternetComputer(DFINITY)"
$ws.Range("C23").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"

# Row -> Price (column D), Volume(1h) (column E)
$updates = @(
    @{ Row = 2;  D = "42.225.03"; E = "  +1.46%  " },
    @{ Row = 3;  D = "2.171.40";  E = "  +0.29%  " },
    @{ Row = 4;  E = "  -0.10%  " },
    @{ Row = 5;  D = "253.39";    E = "  +6.53%  " },
    @{ Row = 6;  D = "0.610";     E = "  +0.15%  " },
    @{ Row = 7;  D = "73.34";     E = "  +1.57%  " },
    @{ Row = 8;  E = "  -0.04%  " },
    @{ Row = 9;  D = "0.580";     E = "  +0.31%  " },
    @{ Row = 10; D = "39.79";     E = "  +0.12%  " },
    @{ Row = 11; D = "0.0907";    E = "  +0.07%  " },
    @{ Row = 12; E = "  +0.63%  " },
    @{ Row = 13; D = "6.74";      E = "  +0.65%  " },
    @{ Row = 14; D = "2.495.70";  E = "  +0.22%  " },
    @{ Row = 15; D = "14.18";     E = "  -0.67%  " },
    @{ Row = 16; D = "2.173.18";  E = "  +0.88%  " },
    @{ Row = 17; D = "0.763";     E = "  -1.72%  " },
    @{ Row = 18; D = "42.098.47"; E = "  +1.43%  " },
    @{ Row = 19; E = "  -0.22%  " },
    @{ Row = 20; D = "70.43";     E = "  +0.63%  " },
    @{ Row = 21; D = "5.82";      E = "  +0.75%  " },
    @{ Row = 22; D = "226.29";    E = "  +0.03%  " },
    @{ Row = 23; D = "9.53";      E = "  -3.91%  " },
    @{ Row = 24; D = "2.14";      E = "  +6.50%  " },
    @{ Row = 25; E = "  -0.17%  " },
    @{ Row = 26; D = "10.43";     E = "  -2.59%  " },
    @{ Row = 27; E = "  +1.44%  " },
    @{ Row = 28; E = "  +2.67%  " },
    @{ Row = 29; D = "2.19";      E = "  +0.19%  " },
    @{ Row = 30; D = "36.57";     E = "  +11.63%  " },
    @{ Row = 31; D = "168.33";    E = "  -1.61%  " },
    @{ Row = 32; D = "19.95";     E = "  +0.70%  " },
    @{ Row = 33; E = "  +4.41%  " },
    @{ Row = 34; D = "5.11";      E = "  -5.14%  " },
    @{ Row = 35; E = "  -0.30%  " },
    @{ Row = 36; D = "0.107";     E = "  +3.92%  " },
    @{ Row = 37; D = "4.23";      E = "  -1.09%  " },
    @{ Row = 38; E = "  +8.85%  " },
    @{ Row = 39; D = "11.81";     E = "  -2.92%  " },
    @{ Row = 40; D = "2.04";      E = "  -2.47%  " },
    @{ Row = 41; E = "  +3.75%  " },
    @{ Row = 42; D = "58.69";     E = "  -0.08%  " },
    @{ Row = 43; D = "5.12";      E = "  -4.52%  " },
    @{ Row = 44; D = "102.13";    E = "  +4.96%  " },
    @{ Row = 45; D = "0.463";     E = "  +15.03%  " },
    @{ Row = 46; D = "8.22";      E = "  -2.53%  " },
    @{ Row = 47; D = "0.0965";    E = "  +0.20%  " },
    @{ Row = 48; D = "2.39";      E = "  +9.54%  " },
    @{ Row = 49; E = "  +0.60%  " },
    @{ Row = 50; E = "  +0.51%  " },
    @{ Row = 51; E = "  +0.98%  " }
)

foreach ($u in $updates) {
    if ($u.ContainsKey("D")) {
        Set-TextCell "D$($u.Row)" $u.D
    }
    if ($u.ContainsKey("E")) {
        $ws.Range("E$($u.Row)").Value = $u.E
    }
}
